$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.640.36'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '1.870.25'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4634'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3884'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07873'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9761'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Value = '1.855.41'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.998'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.704'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06986'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.28%  '
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '28.632.87'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.283'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.116'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '2.123.12'
$ws.Range("E25").Value = '  +3.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.806'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.989'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09335'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9149'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.269'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.339'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.331'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05779'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02107'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.153'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.768'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5630'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.774'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07190'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5312'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.140'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.140'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.832'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.405'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.34%  '
